$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-05-12 Monday" "2025-05-13 Tuesday"

Replace-Text "308÷5=61, 3" "829÷6=138, 1"
Replace-Text "820÷6=136, 4" "603÷6=100, 3"
Replace-Text "731÷4=182, 3" "336÷3=112, 0"
Replace-Text "536÷4=134, 0" "591÷2=295, 1"
Replace-Text "734÷6=122, 2" "849÷4=212, 1"

Replace-Text "794÷9=88, 2" "141÷9=15, 6"
Replace-Text "550÷8=68, 6" "651÷5=130, 1"
Replace-Text "382÷3=127, 1" "506÷8=63, 2"
Replace-Text "463÷4=115, 3" "508÷6=84, 4"
Replace-Text "370÷9=41, 1" "677÷4=169, 1"

Replace-Text "711÷5=142, 1" "489÷5=97, 4"
Replace-Text "295÷7=42, 1" "213÷8=26, 5"
Replace-Text "378÷4=94, 2" "707÷2=353, 1"
Replace-Text "613÷8=76, 5" "679÷6=113, 1"
Replace-Text "126÷6=21, 0" "601÷8=75, 1"

Replace-Text "758÷7=108, 2" "819÷5=163, 4"
Replace-Text "764÷2=382, 0" "146÷9=16, 2"
Replace-Text "118÷4=29, 2" "533÷2=266, 1"
Replace-Text "645÷6=107, 3" "851÷9=94, 5"
Replace-Text "789÷6=131, 3" "547÷2=273, 1"

Replace-Text "465÷8=58, 1" "377÷4=94, 1"
Replace-Text "636÷7=90, 6" "978÷3=326, 0"
Replace-Text "491÷7=70, 1" "814÷9=90, 4"
Replace-Text "826÷3=275, 1" "402÷5=80, 2"
Replace-Text "568÷9=63, 1" "962÷2=481, 0"

Write-Output "Replacements complete"
